# Update "想去人数" (want-to-go count) figures in column F for the sheets
# that contain the full data table: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 131
    4  = 1646
    5  = 635
    7  = 20
    8  = 11622
    9  = 31
    10 = 94
    11 = 460
    13 = 1095
    15 = 12625
    16 = 13184
    17 = 35
    20 = 36
    21 = 254
    24 = 135
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
